$wb = $excel.ActiveWorkbook

# --- settings sheet: form_title / form_id V2 -> V3 ---
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Range("A2").Value = '(2024 Mar) - 3. SCH/STH – Resultats V3'
$wsSettings.Range("B2").Value = 'sn_sch_sth_impact_20403_3_res_v3'

# --- survey sheet: repeat group name + relevant-formula updates ---
$wsSurvey = $wb.Worksheets.Item("survey")
$wsSurvey.Range("B8").Value = 'sn_r_202403_v3'

# H22-H45 relevant formula: Uniquement -> + Kato katz
$wsSurvey.Range("H22").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H23").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H24").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H25").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H26").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H27").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H29").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H30").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H31").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H32").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H33").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H35").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H36").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H37").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H38").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H39").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H41").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H42").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H43").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H44").Value = '${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H45").Value = '${r_test} = ''Filtration d urine + Kato katz'''

# H48-H51 relevant formula: Uniquement -> + Kato katz
$wsSurvey.Range("H48").Value = '${r_autre}=''Oui'' and ${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H49").Value = '${r_autre}=''Oui'' and ${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H50").Value = '${r_autre}=''Oui'' and ${r_test} = ''Filtration d urine + Kato katz'''
$wsSurvey.Range("H51").Value = '${r_autre}=''Oui'' and ${r_test} = ''Filtration d urine + Kato katz'''

# --- view-state: restore active cell / scroll position to match the authors final selection ---
$wsChoices = $wb.Worksheets.Item("choices")
$wsChoices.Range("B13").Select()

$wsSurvey.Activate()
$wsSurvey.Range("B8").Select()

Write-Output "done"
